# "Generate Report for Handoff"
# Updates the localization-status workbook to reflect a new handoff report run:
#  - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#  - Timestamps for the handoff/handback generation are refreshed
#  - zh-cn Priority changes from "ht" (human translation) to "mt" (machine translation)
#  - A new "Error Detail" message is recorded for the 0351e35a... file, noting the
#    handback file used is stale compared to the latest source revision
#  - Column widths for the Status/Error-Detail columns are resized

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Ready for handoff"
$newHoDate = "2016-11-03 20:26:51"
$newPriority = "mt"
$newZhHandoffDate = "2016-11-03 20:26:37"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb9251b23a89679a096d069d7da044ee1344f02e/e2e/0351e35a-635a-48de-af8c-b4ade49b12f0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a315cc56954762fa1e0ce1d2745e1699bc0c7482/e2e/0351e35a-635a-48de-af8c-b4ade49b12f0.md."

# ---- Overview sheet ----
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Range("G2").Value = $newHoDate
$ws1.Range("G3").Value = $newHoDate

# ---- zh-cn sheet ----
$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus
$ws2.Range("E2").Value = $newPriority
$ws2.Range("E3").Value = $newPriority
$ws2.Range("H2").Value = $newZhHandoffDate
$ws2.Range("H3").Value = $newZhHandoffDate
$ws2.Range("P2").Value = $errorDetail

# ---- de-de sheet ----
$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus
$ws3.Range("E2").Value = $newPriority
$ws3.Range("E3").Value = $newPriority
$ws3.Range("H2").Value = $newHoDate
$ws3.Range("H3").Value = $newHoDate
$ws3.Range("P2").Value = $errorDetail

# ---- Column width adjustments ----
# Overview: Status columns (E, F) get narrower
$ws1.Range("E1").ColumnWidth = 16.333333333333336
$ws1.Range("F1").ColumnWidth = 16.333333333333336

# zh-cn / de-de: Status column (C) narrower, Error Detail column (P) wider
$ws2.Range("C1").ColumnWidth = 16.333333333333336
$ws2.Range("P1").ColumnWidth = 39.16666666666667

$ws3.Range("C1").ColumnWidth = 16.333333333333336
$ws3.Range("P1").ColumnWidth = 39.16666666666667
